$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Piscinas" row as an available service by flagging the
# "Desplegar" column (A7) with an X.
$ws.Range("A7").Value = "X"

# Leave the selection on the cell that was just edited.
$ws.Range("A7").Select()
